$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Lancers scrape refreshed at 2025-09-20 12:38:21 JST: two new postings merged in
# ("competitive AI" + "education WEB site"), whole table re-sorted by G (priority
# score) descending, and every row timestamp bumped to the new scrape time.

# Rebuild hyperlinks from scratch: delete existing column-F links so relationship
# ids do not pile up stale entries, then re-add them after the new URLs are written.
$ws.Range("F2:F7").Hyperlinks.Delete()

$rows = @(
    @{ A='2025-09-20 12:38:21'; B='競馬AIの開発ができる方、もしくはすでに開発済みの方'; C='システム開発'; D='300,000 円 ~ 500,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5397293'; G=375; H='🔥AI,Ai ◆開発' },
    @{ A='2025-09-20 12:38:21'; B='【急募】スマホアプリ自動化デモ開発(LLM連携)'; C='システム開発'; D='200,000 円 ~ 300,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5397035'; G=228; H='★スマホアプリ ◆開発,自動化 ◇アプリ' },
    @{ A='2025-09-20 12:38:21'; B='システム開発の案件紹介'; C='システム開発'; D='500,000 円 ~ 1,000,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5397121'; G=125; H='◆開発,システム開発' },
    @{ A='2025-09-20 12:38:21'; B='システム開発において活躍できる案件紹介'; C='システム開発'; D='500,000 円 ~ 1,000,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5397117'; G=125; H='◆開発,システム開発' },
    @{ A='2025-09-20 12:38:21'; B='システム開発の複数案件紹介'; C='システム開発'; D='500,000 円 ~ 1,000,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5397127'; G=125; H='◆開発,システム開発' },
    @{ A='2025-09-20 12:38:21'; B='【急募】教育系のWEBサイトの作成'; C='システム開発'; D='20,000 円 ~ 50,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5397264'; G=33; H='◇サイト' },
    @{ A='2025-09-20 12:38:21'; B='【フォートナイト】クリエイティブ作品を世界に公開したい!'; C='システム開発'; D='50,000 円 ~ 100,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5397192'; G=18; H=$null },
    @{ A='2025-09-20 12:38:21'; B='初回 Web広告のタグ設置・動作確認'; C='システム開発'; D='50,000 円 ~ 100,000 円 / 固定'; E='期限情報なし'; F='https://www.lancers.jp/work/detail/5397007'; G=18; H=$null }
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $d = $rows[$i]
    $ws.Range("A$r").Value = $d.A
    $ws.Range("B$r").Value = $d.B
    $ws.Range("C$r").Value = $d.C
    $ws.Range("D$r").Value = $d.D
    $ws.Range("E$r").Value = $d.E
    $ws.Range("G$r").Value = $d.G
    if ($null -eq $d.H) {
        $ws.Range("H$r").ClearContents()
    } else {
        $ws.Range("H$r").Value = $d.H
    }

    $cellF = $ws.Range("F$r")
    $cellF.Value = $d.F
    $ws.Hyperlinks.Add($cellF, $d.F)
    $cellF.Style = "Hyperlink"
}
